$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LR-pair rows (3 sending clusters x 5 target clusters) with refreshed TPM-based metrics
$data = @(
    @("ECs", "Tff3", "Cxcr4", "ECs", [double]"3", [double]"1", [double]"1.020543333333333", [double]"3.06163", [double]"0.2074480126456944", [double]"0.2074480126456944", [double]"3", [double]"1", [double]"27.681071", [double]"83.04321300000001", [double]"0.05045805550111082", [double]"0.05045805550111081", [double]"28.24973246857667", [double]"254.24759221719", [double]"0.01046742333567159", [double]"0.01046742333567159"),
    @("ECs", "Tff3", "Cxcr4", "FAPs", [double]"3", [double]"1", [double]"1.020543333333333", [double]"3.06163", [double]"0.2074480126456944", [double]"0.2074480126456944", [double]"2", [double]"0.6666666666666666", [double]"0.097204", [double]"0.291612", [double]"0.0001771869602491167", [double]"0.0001771869602491166", [double]"0.09920089417333333", [double]"0.89280804756", [double]"3.675708277041091E-05", [double]"3.67570827704109E-05"),
    @("ECs", "Tff3", "Cxcr4", "Inflammatory-Mac", [double]"3", [double]"1", [double]"1.020543333333333", [double]"3.06163", [double]"0.2074480126456944", [double]"0.2074480126456944", [double]"3", [double]"1", [double]"272.2666776666667", [double]"816.800033", [double]"0.4962975288350554", [double]"0.4962975288350553", [double]"277.8599427815323", [double]"2500.73948503379", [double]"0.1029559360378015", [double]"0.1029559360378015"),
    @("ECs", "Tff3", "Cxcr4", "MuSCs", [double]"3", [double]"1", [double]"1.020543333333333", [double]"3.06163", [double]"0.2074480126456944", [double]"0.2074480126456944", [double]"3", [double]"1", [double]"11.73516533333333", [double]"35.205496", [double]"0.02139128300722342", [double]"0.02139128300722341", [double]"11.97624474649778", [double]"107.78620271848", [double]"0.004437579147790111", [double]"0.004437579147790111"),
    @("ECs", "Tff3", "Cxcr4", "Resolving-Mac", [double]"3", [double]"1", [double]"1.020543333333333", [double]"3.06163", [double]"0.2074480126456944", [double]"0.2074480126456944", [double]"3", [double]"1", [double]"236.8155566666667", [double]"710.44667", [double]"0.4316759456963613", [double]"0.4316759456963613", [double]"241.6805375857889", [double]"2175.1248382721", [double]"0.08955031704166085", [double]"0.08955031704166085"),
    @("Inflammatory-Mac", "Tff3", "Cxcr4", "ECs", [double]"3", [double]"1", [double]"2.383042", [double]"7.149126", [double]"0.4844060127623727", [double]"0.4844060127623726", [double]"3", [double]"1", [double]"27.681071", [double]"83.04321300000001", [double]"0.05045805550111082", [double]"0.05045805550111081", [double]"65.96515479798201", [double]"593.6863931818381", [double]"0.0244421854770356", [double]"0.02444218547703559"),
    @("Inflammatory-Mac", "Tff3", "Cxcr4", "FAPs", [double]"3", [double]"1", [double]"2.383042", [double]"7.149126", [double]"0.4844060127623727", [double]"0.4844060127623726", [double]"2", [double]"0.6666666666666666", [double]"0.097204", [double]"0.291612", [double]"0.0001771869602491167", [double]"0.0001771869602491166", [double]"0.231641214568", [double]"2.084770931112", [double]"8.583042892775961E-05", [double]"8.58304289277596E-05"),
    @("Inflammatory-Mac", "Tff3", "Cxcr4", "Inflammatory-Mac", [double]"3", [double]"1", [double]"2.383042", [double]"7.149126", [double]"0.4844060127623727", [double]"0.4844060127623726", [double]"3", [double]"1", [double]"272.2666776666667", [double]"816.800033", [double]"0.4962975288350554", [double]"0.4962975288350553", [double]"648.8229280801287", [double]"5839.406352721157", [double]"0.2404095070868079", [double]"0.2404095070868078"),
    @("Inflammatory-Mac", "Tff3", "Cxcr4", "MuSCs", [double]"3", [double]"1", [double]"2.383042", [double]"7.149126", [double]"0.4844060127623727", [double]"0.4844060127623726", [double]"3", [double]"1", [double]"11.73516533333333", [double]"35.205496", [double]"0.02139128300722342", [double]"0.02139128300722341", [double]"27.96539186627733", [double]"251.688526796496", [double]"0.01036206610940059", [double]"0.01036206610940059"),
    @("Inflammatory-Mac", "Tff3", "Cxcr4", "Resolving-Mac", [double]"3", [double]"1", [double]"2.383042", [double]"7.149126", [double]"0.4844060127623727", [double]"0.4844060127623726", [double]"3", [double]"1", [double]"236.8155566666667", [double]"710.44667", [double]"0.4316759456963613", [double]"0.4316759456963613", [double]"564.3414177900468", [double]"5079.07276011042", [double]"0.2091064236602009", [double]"0.2091064236602008"),
    @("Resolving-Mac", "Tff3", "Cxcr4", "ECs", [double]"3", [double]"1", [double]"1.515928333333333", [double]"4.547785", [double]"0.3081459745919329", [double]"0.3081459745919329", [double]"3", [double]"1", [double]"27.681071", [double]"83.04321300000001", [double]"0.05045805550111082", [double]"0.05045805550111081", [double]"41.96251982591167", [double]"377.662678433205", [double]"0.01554844668840363", [double]"0.01554844668840363"),
    @("Resolving-Mac", "Tff3", "Cxcr4", "FAPs", [double]"3", [double]"1", [double]"1.515928333333333", [double]"4.547785", [double]"0.3081459745919329", [double]"0.3081459745919329", [double]"2", [double]"0.6666666666666666", [double]"0.097204", [double]"0.291612", [double]"0.0001771869602491167", [double]"0.0001771869602491166", [double]"0.1473542977133333", [double]"1.32618867942", [double]"5.459944855094612E-05", [double]"5.459944855094611E-05"),
    @("Resolving-Mac", "Tff3", "Cxcr4", "Inflammatory-Mac", [double]"3", [double]"1", [double]"1.515928333333333", [double]"4.547785", [double]"0.3081459745919329", [double]"0.3081459745919329", [double]"3", [double]"1", [double]"272.2666776666667", [double]"816.800033", [double]"0.4962975288350554", [double]"0.4962975288350553", [double]"412.7367708974339", [double]"3714.630938076905", [double]"0.1529320857104461", [double]"0.152932085710446"),
    @("Resolving-Mac", "Tff3", "Cxcr4", "MuSCs", [double]"3", [double]"1", [double]"1.515928333333333", [double]"4.547785", [double]"0.3081459745919329", [double]"0.3081459745919329", [double]"3", [double]"1", [double]"11.73516533333333", [double]"35.205496", [double]"0.02139128300722342", [double]"0.02139128300722341", [double]"17.78966962515111", [double]"160.10702662636", [double]"0.006591637750032712", [double]"0.006591637750032711"),
    @("Resolving-Mac", "Tff3", "Cxcr4", "Resolving-Mac", [double]"3", [double]"1", [double]"1.515928333333333", [double]"4.547785", [double]"0.3081459745919329", [double]"0.3081459745919329", [double]"3", [double]"1", [double]"236.8155566666667", [double]"710.44667", [double]"0.4316759456963613", [double]"0.4316759456963613", [double]"358.9954121251056", [double]"3230.95870912595", [double]"0.1330192049944995", [double]"0.1330192049944995"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

